$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2026-02-15 Sunday" "2026-02-16 Monday"

Replace-Text "416÷2=" "821÷4="
Replace-Text "445÷3=" "225÷2="
Replace-Text "186÷5=" "252÷6="
Replace-Text "287÷2=" "900÷3="
Replace-Text "360÷3=" "963÷7="

Replace-Text "437÷3=" "318÷5="
Replace-Text "599÷2=" "342÷5="
Replace-Text "964÷4=" "255÷2="
Replace-Text "554÷8=" "378÷5="
Replace-Text "831÷3=" "725÷7="

Replace-Text "671÷2=" "956÷3="
Replace-Text "705÷4=" "460÷3="
Replace-Text "188÷2=" "604÷2="
Replace-Text "431÷3=" "983÷6="
Replace-Text "935÷7=" "499÷5="

Replace-Text "334÷8=" "972÷4="
Replace-Text "419÷4=" "544÷6="
Replace-Text "662÷9=" "871÷6="
Replace-Text "720÷3=" "805÷6="
Replace-Text "147÷6=" "162÷3="

Replace-Text "321÷3=" "549÷9="
Replace-Text "164÷3=" "878÷3="
Replace-Text "418÷3=" "900÷4="
Replace-Text "568÷9=" "624÷9="
Replace-Text "717÷4=" "160÷8="
